$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$r = $t.Cell(1,1).Range
$tr = $d.Range($r.Start, $r.End - 1)
$tr.Text = "55-12="
$r = $t.Cell(1,2).Range
$tr = $d.Range($r.Start, $r.End - 1)
$tr.Text = "82+3="
$r = $t.Cell(1,3).Range
$tr = $d.Range($r.Start, $r.End - 1)
$tr.Text = "37+51="
$r = $t.Cell(1,4).Range
$tr = $d.Range($r.Start, $r.End - 1)
$tr.Text = "58+37="
$r = $t.Cell(1,5).Range
$tr = $d.Range($r.Start, $r.End - 1)
$tr.Text = "10+45="

$r = $t.Cell(2,1).Range
$tr = $d.Range($r.Start, $r.End - 1)
$tr.Text = "40+20="
$r = $t.Cell(2,2).Range
$tr = $d.Range($r.Start, $r.End - 1)
$tr.Text = "17+33="
$r = $t.Cell(2,3).Range
$tr = $d.Range($r.Start, $r.End - 1)
$tr.Text = "60-47="
$r = $t.Cell(2,4).Range
$tr = $d.Range($r.Start, $r.End - 1)
$tr.Text = "5+86="
$r = $t.Cell(2,5).Range
$tr = $d.Range($r.Start, $r.End - 1)
$tr.Text = "51-21="

$r = $t.Cell(3,1).Range
$tr = $d.Range($r.Start, $r.End - 1)
$tr.Text = "96-96="
$r = $t.Cell(3,2).Range
$tr = $d.Range($r.Start, $r.End - 1)
$tr.Text = "79-28="
$r = $t.Cell(3,3).Range
$tr = $d.Range($r.Start, $r.End - 1)
$tr.Text = "50-16="
$r = $t.Cell(3,4).Range
$tr = $d.Range($r.Start, $r.End - 1)
$tr.Text = "77+13="
$r = $t.Cell(3,5).Range
$tr = $d.Range($r.Start, $r.End - 1)
$tr.Text = "57-5="

$r = $t.Cell(4,1).Range
$tr = $d.Range($r.Start, $r.End - 1)
$tr.Text = "61-48="
$r = $t.Cell(4,2).Range
$tr = $d.Range($r.Start, $r.End - 1)
$tr.Text = "89-74="
$r = $t.Cell(4,3).Range
$tr = $d.Range($r.Start, $r.End - 1)
$tr.Text = "62-29="
$r = $t.Cell(4,4).Range
$tr = $d.Range($r.Start, $r.End - 1)
$tr.Text = "19+57="
$r = $t.Cell(4,5).Range
$tr = $d.Range($r.Start, $r.End - 1)
$tr.Text = "55-42="

$r = $t.Cell(5,1).Range
$tr = $d.Range($r.Start, $r.End - 1)
$tr.Text = "11+87="
$r = $t.Cell(5,2).Range
$tr = $d.Range($r.Start, $r.End - 1)
$tr.Text = "83-8="
$r = $t.Cell(5,3).Range
$tr = $d.Range($r.Start, $r.End - 1)
$tr.Text = "21+32="
$r = $t.Cell(5,4).Range
$tr = $d.Range($r.Start, $r.End - 1)
$tr.Text = "55+14="
$r = $t.Cell(5,5).Range
$tr = $d.Range($r.Start, $r.End - 1)
$tr.Text = "86-63="

$r = $t.Cell(6,1).Range
$tr = $d.Range($r.Start, $r.End - 1)
$tr.Text = "58+20="
$r = $t.Cell(6,2).Range
$tr = $d.Range($r.Start, $r.End - 1)
$tr.Text = "83-26="
$r = $t.Cell(6,3).Range
$tr = $d.Range($r.Start, $r.End - 1)
$tr.Text = "54-7="
$r = $t.Cell(6,4).Range
$tr = $d.Range($r.Start, $r.End - 1)
$tr.Text = "44+10="
$r = $t.Cell(6,5).Range
$tr = $d.Range($r.Start, $r.End - 1)
$tr.Text = "20+35="

$r = $t.Cell(7,1).Range
$tr = $d.Range($r.Start, $r.End - 1)
$tr.Text = "16+82="
$r = $t.Cell(7,2).Range
$tr = $d.Range($r.Start, $r.End - 1)
$tr.Text = "81-69="
$r = $t.Cell(7,3).Range
$tr = $d.Range($r.Start, $r.End - 1)
$tr.Text = "30+33="
$r = $t.Cell(7,4).Range
$tr = $d.Range($r.Start, $r.End - 1)
$tr.Text = "53-17="
$r = $t.Cell(7,5).Range
$tr = $d.Range($r.Start, $r.End - 1)
$tr.Text = "44-12="

$r = $t.Cell(8,1).Range
$tr = $d.Range($r.Start, $r.End - 1)
$tr.Text = "33+53="
$r = $t.Cell(8,2).Range
$tr = $d.Range($r.Start, $r.End - 1)
$tr.Text = "51+44="
$r = $t.Cell(8,3).Range
$tr = $d.Range($r.Start, $r.End - 1)
$tr.Text = "16-16="
$r = $t.Cell(8,4).Range
$tr = $d.Range($r.Start, $r.End - 1)
$tr.Text = "35+28="
$r = $t.Cell(8,5).Range
$tr = $d.Range($r.Start, $r.End - 1)
$tr.Text = "92-82="

$r = $t.Cell(9,1).Range
$tr = $d.Range($r.Start, $r.End - 1)
$tr.Text = "78-60="
$r = $t.Cell(9,2).Range
$tr = $d.Range($r.Start, $r.End - 1)
$tr.Text = "70+10="
$r = $t.Cell(9,3).Range
$tr = $d.Range($r.Start, $r.End - 1)
$tr.Text = "99-97="
$r = $t.Cell(9,4).Range
$tr = $d.Range($r.Start, $r.End - 1)
$tr.Text = "30-0="
$r = $t.Cell(9,5).Range
$tr = $d.Range($r.Start, $r.End - 1)
$tr.Text = "10+9="

$r = $t.Cell(10,1).Range
$tr = $d.Range($r.Start, $r.End - 1)
$tr.Text = "36+17="
$r = $t.Cell(10,2).Range
$tr = $d.Range($r.Start, $r.End - 1)
$tr.Text = "52-47="
$r = $t.Cell(10,3).Range
$tr = $d.Range($r.Start, $r.End - 1)
$tr.Text = "32+22="
$r = $t.Cell(10,4).Range
$tr = $d.Range($r.Start, $r.End - 1)
$tr.Text = "69-16="
$r = $t.Cell(10,5).Range
$tr = $d.Range($r.Start, $r.End - 1)
$tr.Text = "80-17="

$r = $t.Cell(11,1).Range
$tr = $d.Range($r.Start, $r.End - 1)
$tr.Text = "63+8="
$r = $t.Cell(11,2).Range
$tr = $d.Range($r.Start, $r.End - 1)
$tr.Text = "59-18="
$r = $t.Cell(11,3).Range
$tr = $d.Range($r.Start, $r.End - 1)
$tr.Text = "60-25="
$r = $t.Cell(11,4).Range
$tr = $d.Range($r.Start, $r.End - 1)
$tr.Text = "36+48="
$r = $t.Cell(11,5).Range
$tr = $d.Range($r.Start, $r.End - 1)
$tr.Text = "28+44="

$r = $t.Cell(12,1).Range
$tr = $d.Range($r.Start, $r.End - 1)
$tr.Text = "72-54="
$r = $t.Cell(12,2).Range
$tr = $d.Range($r.Start, $r.End - 1)
$tr.Text = "48+21="
$r = $t.Cell(12,3).Range
$tr = $d.Range($r.Start, $r.End - 1)
$tr.Text = "81-18="
$r = $t.Cell(12,4).Range
$tr = $d.Range($r.Start, $r.End - 1)
$tr.Text = "83-6="
$r = $t.Cell(12,5).Range
$tr = $d.Range($r.Start, $r.End - 1)
$tr.Text = "12+64="

$r = $t.Cell(13,1).Range
$tr = $d.Range($r.Start, $r.End - 1)
$tr.Text = "65-11="
$r = $t.Cell(13,2).Range
$tr = $d.Range($r.Start, $r.End - 1)
$tr.Text = "75-58="
$r = $t.Cell(13,3).Range
$tr = $d.Range($r.Start, $r.End - 1)
$tr.Text = "98-31="
$r = $t.Cell(13,4).Range
$tr = $d.Range($r.Start, $r.End - 1)
$tr.Text = "4+16="
$r = $t.Cell(13,5).Range
$tr = $d.Range($r.Start, $r.End - 1)
$tr.Text = "41+48="

$r = $t.Cell(14,1).Range
$tr = $d.Range($r.Start, $r.End - 1)
$tr.Text = "45+31="
$r = $t.Cell(14,2).Range
$tr = $d.Range($r.Start, $r.End - 1)
$tr.Text = "4+51="
$r = $t.Cell(14,3).Range
$tr = $d.Range($r.Start, $r.End - 1)
$tr.Text = "48-37="
$r = $t.Cell(14,4).Range
$tr = $d.Range($r.Start, $r.End - 1)
$tr.Text = "46+45="
$r = $t.Cell(14,5).Range
$tr = $d.Range($r.Start, $r.End - 1)
$tr.Text = "99-98="

$r = $t.Cell(15,1).Range
$tr = $d.Range($r.Start, $r.End - 1)
$tr.Text = "66+4="
$r = $t.Cell(15,2).Range
$tr = $d.Range($r.Start, $r.End - 1)
$tr.Text = "67+26="
$r = $t.Cell(15,3).Range
$tr = $d.Range($r.Start, $r.End - 1)
$tr.Text = "70-68="
$r = $t.Cell(15,4).Range
$tr = $d.Range($r.Start, $r.End - 1)
$tr.Text = "33+14="
$r = $t.Cell(15,5).Range
$tr = $d.Range($r.Start, $r.End - 1)
$tr.Text = "55-49="

$r = $t.Cell(16,1).Range
$tr = $d.Range($r.Start, $r.End - 1)
$tr.Text = "63-17="
$r = $t.Cell(16,2).Range
$tr = $d.Range($r.Start, $r.End - 1)
$tr.Text = "57-13="
$r = $t.Cell(16,3).Range
$tr = $d.Range($r.Start, $r.End - 1)
$tr.Text = "3+20="
$r = $t.Cell(16,4).Range
$tr = $d.Range($r.Start, $r.End - 1)
$tr.Text = "93-31="
$r = $t.Cell(16,5).Range
$tr = $d.Range($r.Start, $r.End - 1)
$tr.Text = "76-70="

$r = $t.Cell(17,1).Range
$tr = $d.Range($r.Start, $r.End - 1)
$tr.Text = "97-5="
$r = $t.Cell(17,2).Range
$tr = $d.Range($r.Start, $r.End - 1)
$tr.Text = "61-0="
$r = $t.Cell(17,3).Range
$tr = $d.Range($r.Start, $r.End - 1)
$tr.Text = "56-42="
$r = $t.Cell(17,4).Range
$tr = $d.Range($r.Start, $r.End - 1)
$tr.Text = "8+50="
$r = $t.Cell(17,5).Range
$tr = $d.Range($r.Start, $r.End - 1)
$tr.Text = "29+22="

$r = $t.Cell(18,1).Range
$tr = $d.Range($r.Start, $r.End - 1)
$tr.Text = "25-20="
$r = $t.Cell(18,2).Range
$tr = $d.Range($r.Start, $r.End - 1)
$tr.Text = "61-54="
$r = $t.Cell(18,3).Range
$tr = $d.Range($r.Start, $r.End - 1)
$tr.Text = "76-12="
$r = $t.Cell(18,4).Range
$tr = $d.Range($r.Start, $r.End - 1)
$tr.Text = "0+56="
$r = $t.Cell(18,5).Range
$tr = $d.Range($r.Start, $r.End - 1)
$tr.Text = "27+51="

$r = $t.Cell(19,1).Range
$tr = $d.Range($r.Start, $r.End - 1)
$tr.Text = "39-8="
$r = $t.Cell(19,2).Range
$tr = $d.Range($r.Start, $r.End - 1)
$tr.Text = "51-14="
$r = $t.Cell(19,3).Range
$tr = $d.Range($r.Start, $r.End - 1)
$tr.Text = "25-0="
$r = $t.Cell(19,4).Range
$tr = $d.Range($r.Start, $r.End - 1)
$tr.Text = "60-19="
$r = $t.Cell(19,5).Range
$tr = $d.Range($r.Start, $r.End - 1)
$tr.Text = "36+34="

$r = $t.Cell(20,1).Range
$tr = $d.Range($r.Start, $r.End - 1)
$tr.Text = "27+71="
$r = $t.Cell(20,2).Range
$tr = $d.Range($r.Start, $r.End - 1)
$tr.Text = "59-23="
$r = $t.Cell(20,3).Range
$tr = $d.Range($r.Start, $r.End - 1)
$tr.Text = "12+81="
$r = $t.Cell(20,4).Range
$tr = $d.Range($r.Start, $r.End - 1)
$tr.Text = "71-31="
$r = $t.Cell(20,5).Range
$tr = $d.Range($r.Start, $r.End - 1)
$tr.Text = "13+15="
